$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "PIN" column values (text like "4 (PD2)") with plain numeric values.
$ws.Range("D2").Value = 4
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 6

# Update the selected cell/range shown in the sheet view.
$ws.Range("D6").Select()
